# Starting to call the sample metadata update service: replace the
# sample_id column with sample names pulled back from the service, clear
# out the (now stale) third data row, and leave the selection sitting on
# that cleared row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: sample_id -> sample_name
$ws.Range("A1").Value = "sample_name"

# Replace the numeric sample ids with the names returned by the service
$ws.Range("A2").Value = "Project 1 Sample 1"
$ws.Range("A3").Value = "Project 1 Sample 2"

# Row 4 had no corresponding sample from the service - clear it out
$ws.Range("A4:D4").ClearContents()

# Touch the formatting on the whole used range (mirrors re-applying the
# "Normal" style across the table once the update finished)
$ws.Range("A1:D4").Style = "Normal"

# Leave the selection on the now-empty row
[void]$ws.Range("A4").Select()
